$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- [1.17.2] ~ separate the logic for printing envelops from printing postcards
#              ~ use the postcards templates with white background --------------
# Insert the new version row right after the last existing entry (row 16).
$ws.Rows.Item(17).Insert(-4121)

$ws.Range("A17").Value = "[1.17.2]"
$ws.Range("B17").Value = "~ separate the logic for printing envelops from printing postcards`n~ use the postcards templates with white background"

# Give the Date cell the same number format / alignment as the other date
# entries in the table (left/top aligned, d-mmm-yy).
$ws.Range("C14").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = Get-Date -Year 2018 -Month 8 -Day 24 -Hour 0 -Minute 0 -Second 0
$ws.Rows.Item(17).RowHeight = 30

# Append the trailing blank row that follows the new entry.
$ws.Rows.Item(18).Insert(-4121)
$ws.Range("A5").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("A18:C18").ClearContents()

# Grow the table (and its autofilter) so the two new rows are included.
$tbl = $ws.ListObjects.Item("Table2")
$tbl.Resize($ws.Range("A1:C18"))

# Restore the selection Excel leaves behind after this kind of edit.
$ws.Range("G13").Select()
